$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.56"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.249"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05789"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.501"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.145"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8180"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8494"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1362"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06941"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03152"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02879"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09401"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.749"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001513"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04718"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005960"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006269"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001235"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004615"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00006902"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.514"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.123"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1346"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03652"

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "BKEXToken"

$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1057"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "CEJI"

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002751"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "KickToken"

$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003015"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008064"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005264"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002335"
